$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B34").Value = "test"
